# Fruta / hortaliza, semanal
# Insert a new weekly record for "Naranja - Fukumoto" before the former
# row 529, shifting every following row down by one (dimension grows
# from A1:T560 to A1:T561).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(529).Insert()

$ws.Range("A529").Value = 5
$ws.Range("B529").Value = 'Macroferia Regional de Talca'
$ws.Range("C529").Value = 'Maule'
$ws.Range("D529").Value = 44753
$ws.Range("E529").Value = 7
$ws.Range("F529").Value = 'Fruta'
$ws.Range("G529").Value = 100102
$ws.Range("H529").Value = 'Cítricos'
$ws.Range("I529").Value = 100102005
$ws.Range("J529").Value = 'Naranja'
$ws.Range("K529").Value = 'Fukumoto'
$ws.Range("L529").Value = 'Primera'
$ws.Range("M529").Value = 500
$ws.Range("N529").Value = 7000
$ws.Range("O529").Value = 7000
$ws.Range("P529").Value = 7000
$ws.Range("Q529").Value = '$/bandeja 15 kilos granel'
$ws.Range("R529").Value = 'Provincia de Melipilla'
$ws.Range("S529").Value = 467
$ws.Range("T529").Value = 15
